# Add a "clarifications needed" column to the Details sheet's table (Table2),
# capturing review questions/clarifications against several requirement rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$tbl = $ws.ListObjects.Item(1)

# Extend the table by one column (A2:M13 -> A2:N13)
$newCol = $tbl.ListColumns.Add()

# Header
$ws.Range("N2").Value = "clarifications needed"

# Clarification text for rows 3-9
$ws.Range("N3").Value = "1. what is the use case of locking specific Biometric auths and combinations?`n2. What is the data taken as input from the Individual?`n3. Is there a mechanism to lock OTP Authentication?`n4. need more clarity on a2`n5. Can Lock/Unlock will be perfomed only by OTP authentication of Mobile number or can it be done by email as well?`n6. Need to check with IDA on the process of authentication done, so that the gaps could be covered`n"

$ws.Range("N4").Value = "1. What is the use case of reprinting? Will there be a cost associated to it? If not can it be abused by the individual?will there be a limit on number times an individual can access it? If cost is associated, will there be a check performed for the payment?`n2. Why RID is accepted as an input parameter? What is the use case."

$ws.Range("N5").Value = "1. if demo auth is locked? What happens?`n2. what is use case of providing a RID PDF, why not just a RID number?`n3. why do we have a size check here? Shouldn’t it be stopped at the initial stage?"

$ws.Range("N6").Value = "1. Is this requirement still part of Resident services?`n2. if demo auth is locked? What happens?`n3. why do we have a size check here? Shouldn’t it be stopped at the initial stage?"

$ws.Range("N7").Value = "1. Is this requirement a subset of what registration client does for update? If so can there be an reuse of the feature?`n2. in future if there can be many parameters which can be provisioned for updated? What can be done and is Resident services capable of doing for all? check for scalability?`n"

$ws.Range("N8").Value = "1. why is RID an input here and not UIN?"

$ws.Range("N9").Value = "1.what is security code?"

# Rows 10-13 stay blank but remain part of the table column.

# Formatting: align with the look of the rest of the table (wrapped text,
# top/left aligned, thin border all around) and give the column a sensible width.
$colRange = $ws.Range("N2:N13")
$colRange.WrapText = $true
$colRange.HorizontalAlignment = -4131
$colRange.VerticalAlignment = -4160
$colRange.Borders.LineStyle = 1
$colRange.Borders.Weight = 2

$ws.Columns.Item(14).ColumnWidth = 48.5
